# JS-Frameworks-Self-Evaluation-Protocol.xlsx — "Max-min length and styles added"
#
# - Numbers of Commits in GitHub comment bumped from "49 commits" to "50 commits"
# - AngularJS Project Structure score adjusted from 25 to 23 (Total Score formula
#   in C44 recalculates automatically from 347 to 345)
# - Selection/view moved from C43 (scrolled to row 34) back to E9 at the top of
#   the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Self-Evaluation-Protocol")

# Numbers of Commits in GitHub -> comment cell E9
$ws.Range("E9").Value = "50 commits"

# AngularJS Project Structure score, row 12 (Score column C)
$ws.Range("C12").Value = 23

# Bring the view back to the top of the sheet and leave the selection on E9,
# matching the saved sheetView/selection in the workbook
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("E9").Select()
